$wb = $excel.ActiveWorkbook

# Add a new worksheet "Sheet1" after the last existing sheet (so it lands at the end)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Sheet1"

# ---- Block 1: rows 9-16 ----
$ws.Range("D9").Value = 1000
$ws.Range("E9").Value = 2000
$ws.Range("F9").Value = 5000
$ws.Range("G9").Value = 6500
$ws.Range("H9").Value = 6500
$ws.Range("I9").Value = 6500
$ws.Range("J9").Value = 7000
$ws.Range("K9").Value = 8500
$ws.Range("L9").Value = 9000
$ws.Range("M9").Value = 10500
$ws.Range("N9").Value = 150000000
$ws.Range("D9:N9").Style = "Comma"

$ws.Range("B10").Value = "mean"
$ws.Range("C10").Style = "Comma"
$ws.Range("C10").Formula = "=AVERAGE(D9:N9)"

$ws.Range("B11").Value = "median"
$ws.Range("C11").Style = "Comma"
$ws.Range("C11").Formula = "=MEDIAN(D9:N9)"

$ws.Range("B12").Value = "std"
$ws.Range("C12").Style = "Comma"
$ws.Range("C12").Formula = "=STDEV.S((D9:N9))"

$ws.Range("B13").Value = "q1"
$ws.Range("C13").Style = "Comma"
$ws.Range("C13").Formula = "=QUARTILE.INC(D9:N9,1)"

$ws.Range("B14").Value = "q2"
$ws.Range("C14").Style = "Comma"
$ws.Range("C14").Formula = "=QUARTILE.INC(D9:N9,2)"

$ws.Range("B15").Value = "q3"
$ws.Range("C15").Style = "Comma"
$ws.Range("C15").Formula = "=QUARTILE.INC(D9:N9,3)"

$ws.Range("B16").Value = "modus"
$ws.Range("C16").Style = "Comma"
$ws.Range("C16").Formula = "=MODE.SNGL(D9:N9)"

# ---- Block 2: rows 19-26 ----
$ws.Range("D19").Value = 1000
$ws.Range("E19").Value = 2000
$ws.Range("F19").Value = 5000
$ws.Range("G19").Value = 6500
$ws.Range("H19").Value = 6500
$ws.Range("I19").Value = 6500
$ws.Range("J19").Value = 7000
$ws.Range("K19").Value = 8500
$ws.Range("L19").Value = 9000
$ws.Range("M19").Value = 10500
$ws.Range("N19").Value = 15000
$ws.Range("D19:N19").Style = "Comma"

$ws.Range("B20").Value = "mean"
$ws.Range("C20").Style = "Comma"
$ws.Range("C20").Formula = "=AVERAGE(D19:N19)"

$ws.Range("B21").Value = "median"
$ws.Range("C21").Style = "Comma"
$ws.Range("C21").Formula = "=MEDIAN(D19:N19)"

$ws.Range("B22").Value = "std"
$ws.Range("C22").Style = "Comma"
$ws.Range("C22").Formula = "=STDEV.S((D19:N19))"

$ws.Range("B23").Value = "q1"
$ws.Range("C23").Style = "Comma"
$ws.Range("C23").Formula = "=QUARTILE.INC(D19:N19,1)"

$ws.Range("B24").Value = "q2"
$ws.Range("C24").Style = "Comma"
$ws.Range("C24").Formula = "=QUARTILE.INC(D19:N19,2)"

$ws.Range("B25").Value = "q3"
$ws.Range("C25").Style = "Comma"
$ws.Range("C25").Formula = "=QUARTILE.INC(D19:N19,3)"

$ws.Range("B26").Value = "modus"
$ws.Range("C26").Style = "Comma"
$ws.Range("C26").Formula = "=MODE.SNGL(D19:N19)"

# Column widths matching the target sheet (best-fit style widths, as close
# as the engine's pixel-quantized ColumnWidth grid allows)
$ws.Columns.Item(3).ColumnWidth = 13.5
$ws.Range("D1:L1").EntireColumn.ColumnWidth = 8.666666666666666
$ws.Columns.Item(13).ColumnWidth = 9.666666666666666
$ws.Columns.Item(14).ColumnWidth = 14.5

# Select N20 as the active cell on the new sheet, matching the saved selection
$ws.Range("N20").Select()
